$d = $word.ActiveDocument

# Update the title date line (first paragraph).
$d.Content.Find.Execute("2025-01-02 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-01-03 Friday", 2)

# New arithmetic expressions, in row-major order matching the table's
# existing cell order (20 rows x 5 columns).
$newValues = @(
    "86-47=",
    "43-16=",
    "43+39=",
    "92-48=",
    "82-5=",
    "36+49=",
    "90-89=",
    "29+59=",
    "61-42=",
    "34+9=",
    "18+36=",
    "24+67=",
    "66-7=",
    "7+27=",
    "68+29=",
    "66+7=",
    "49+38=",
    "16+58=",
    "77+17=",
    "68+29=",
    "83-36=",
    "87-9=",
    "15+59=",
    "67+5=",
    "90-66=",
    "9+15=",
    "61-8=",
    "26+66=",
    "87-29=",
    "22-6=",
    "56-48=",
    "6+8=",
    "46+38=",
    "24+69=",
    "36-28=",
    "40-19=",
    "53+18=",
    "91-18=",
    "85-29=",
    "75-56=",
    "65-28=",
    "80-71=",
    "57+26=",
    "30-17=",
    "26-7=",
    "35+38=",
    "95-49=",
    "81-18=",
    "6+5=",
    "32-24=",
    "81-59=",
    "8+46=",
    "51-34=",
    "5+39=",
    "66-8=",
    "34-19=",
    "35+17=",
    "37+54=",
    "29+18=",
    "29+53=",
    "46+35=",
    "85-67=",
    "56+39=",
    "78+17=",
    "66-48=",
    "9+28=",
    "90-17=",
    "83-77=",
    "85-36=",
    "35+38=",
    "5+37=",
    "64+9=",
    "77+4=",
    "9+44=",
    "85-37=",
    "76+5=",
    "34+47=",
    "52+29=",
    "51-22=",
    "37-18=",
    "16+39=",
    "63-18=",
    "40-38=",
    "52-19=",
    "9+35=",
    "42-8=",
    "17+58=",
    "8+18=",
    "28+14=",
    "81-32=",
    "7+15=",
    "6+66=",
    "18+64=",
    "6+78=",
    "15+36=",
    "28+54=",
    "49+9=",
    "39+19=",
    "73-8=",
    "93-59="
)

$tbl = $d.Tables.Item(1)
$rowCount = $tbl.Rows.Count
$colCount = $tbl.Columns.Count

$i = 0
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i++
    }
}
